$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write cells in the same order the shared strings were first introduced:
# A120, A121, D121, C121, C120, D120
$ws.Cells.Item(120, 1).Value = "em_ui_filter"
$ws.Cells.Item(121, 1).Value = "em_ui_add"
$ws.Cells.Item(121, 4).Value = "Add"
$ws.Cells.Item(121, 3).Value = "追加 "
$ws.Cells.Item(120, 3).Value = "最近の会話フィルター"
$ws.Cells.Item(120, 4).Value = "Recent Actions Filter"

$ws.Range("D123").Select()
